$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: H1/I1 were missing the centered header style that the
# rest of row 1 (B1,C1,D1,E1,F1,G1) already has; bring them in line. ---
$ws.Range("H1:I1").HorizontalAlignment = -4108
$ws.Range("H1:I1").VerticalAlignment = -4108

# --- Row 2: fill in the first apartment listing. The "Link" (I2) and
# "Suites" (G2) columns are intentionally left blank - "falta so o link". ---
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "CLNW 2/3, NOROESTE, BRASILIA"
$ws.Range("C2").Value = "R$ 4.300"
$ws.Range("D2").Value = "R$ 81"
$ws.Range("E2").Value = "53 m²"
$ws.Range("F2").Value = "2 Quartos"
$ws.Range("H2").Value = "1 Vaga"

# Right-align + vertically center the newly written text columns (E,F,H)
# to match the sheet's existing "value" style used in C2/D2.
$ws.Range("E2").HorizontalAlignment = -4152
$ws.Range("E2").VerticalAlignment = -4108
$ws.Range("F2").HorizontalAlignment = -4152
$ws.Range("F2").VerticalAlignment = -4108
$ws.Range("H2").HorizontalAlignment = -4152
$ws.Range("H2").VerticalAlignment = -4108

# --- View state: active cell moves to I2 (the still-empty Link cell) ---
$ws.Range("I2").Select()
